$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 53.88
$ws.Range("C3").Value = 146.99
$ws.Range("C4").Value = 199.65
$ws.Range("C5").Value = 439.3
$ws.Range("C6").Value = 408.42
$ws.Range("C7").Value = 149.24
$ws.Range("C8").Value = 818.1799999999999
$ws.Range("C9").Value = 311.63
